$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells keep their text representation
# (values like "1.002" would otherwise be parsed as numbers)
$priceCells = @("D2", "D3", "D5", "D6", "D7", "D8", "D10", "D11", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D45", "D46", "D48", "D49", "D50", "D51")
foreach ($ref in $priceCells) {
    $ws.Range($ref).NumberFormat = "@"
}

# Apply updated values
$ws.Range("D2").Value = '27.408.91'
$ws.Range("E2").Value = '  -3.04%  '
$ws.Range("D3").Value = '1.741.77'
$ws.Range("E3").Value = '  -3.49%  '
$ws.Range("D5").Value = '322.20'
$ws.Range("E5").Value = '  -4.76%  '
$ws.Range("D6").Value = '1.002'
$ws.Range("E6").Value = '  +0.23%  '
$ws.Range("D7").Value = '0.4223'
$ws.Range("E7").Value = '  -10.10%  '
$ws.Range("D8").Value = '0.3574'
$ws.Range("E8").Value = '  -5.79%  '
$ws.Range("E9").Value = '  +0.11%  '
$ws.Range("D10").Value = '0.07393'
$ws.Range("E10").Value = '  -2.98%  '
$ws.Range("D11").Value = '1.109'
$ws.Range("E11").Value = '  -3.68%  '
$ws.Range("E12").Value = '  +0.15%  '
$ws.Range("D13").Value = '21.33'
$ws.Range("E13").Value = '  -4.80%  '
$ws.Range("D14").Value = '6.082'
$ws.Range("E14").Value = '  -3.92%  '
$ws.Range("D15").Value = '7.162'
$ws.Range("E15").Value = '  -3.91%  '
$ws.Range("D16").Value = '1.742.52'
$ws.Range("E16").Value = '  -3.33%  '
$ws.Range("D17").Value = '0.00001062'
$ws.Range("E17").Value = '  -2.84%  '
$ws.Range("D18").Value = '87.18'
$ws.Range("E18").Value = '  +6.56%  '
$ws.Range("D19").Value = '0.06183'
$ws.Range("E19").Value = '  -7.99%  '
$ws.Range("E20").Value = '  +0.18%  '
$ws.Range("D21").Value = '16.79'
$ws.Range("E21").Value = '  -3.60%  '
$ws.Range("D22").Value = '6.082'
$ws.Range("E22").Value = '  -5.15%  '
$ws.Range("D23").Value = '0.5257'
$ws.Range("E23").Value = '  -5.00%  '
$ws.Range("D24").Value = '27.459.96'
$ws.Range("E24").Value = '  -2.85%  '
$ws.Range("D25").Value = '11.58'
$ws.Range("E25").Value = '  -2.38%  '
$ws.Range("D26").Value = '2.321'
$ws.Range("E26").Value = '  -3.56%  '
$ws.Range("D27").Value = '20.32'
$ws.Range("E27").Value = '  -1.98%  '
$ws.Range("D28").Value = '151.75'
$ws.Range("E28").Value = '  -1.44%  '
$ws.Range("D29").Value = '2.352'
$ws.Range("E29").Value = '  -0.79%  '
$ws.Range("D30").Value = '1.940.55'
$ws.Range("E30").Value = '  -3.39%  '
$ws.Range("D31").Value = '125.96'
$ws.Range("E31").Value = '  -5.43%  '
$ws.Range("D32").Value = '1.202'
$ws.Range("E32").Value = '  -4.16%  '
$ws.Range("D33").Value = '5.657'
$ws.Range("E33").Value = '  -3.45%  '
$ws.Range("B34").Value = 'Stellar'
$ws.Range("C34").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D34").Value = '0.09113'
$ws.Range("E34").Value = '  -5.37%  '
$ws.Range("B35").Value = 'HuobiToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D35").Value = '3.679'
$ws.Range("E35").Value = '  -8.84%  '
$ws.Range("D36").Value = '12.58'
$ws.Range("E36").Value = '  +3.77%  '
$ws.Range("D37").Value = '0.02281'
$ws.Range("E37").Value = '  -3.17%  '
$ws.Range("D38").Value = '0.2125'
$ws.Range("E38").Value = '  -5.40%  '
$ws.Range("D39").Value = '5.077'
$ws.Range("E39").Value = '  -3.32%  '
$ws.Range("D40").Value = '0.06068'
$ws.Range("E40").Value = '  -4.77%  '
$ws.Range("D41").Value = '0.6384'
$ws.Range("E41").Value = '  -3.74%  '
$ws.Range("D42").Value = '1.193'
$ws.Range("E42").Value = '  -3.52%  '
$ws.Range("D43").Value = '1.420'
$ws.Range("E43").Value = '  -5.54%  '
$ws.Range("E44").Value = '  +0.20%  '
$ws.Range("D45").Value = '7.863'
$ws.Range("E45").Value = '  -4.64%  '
$ws.Range("D46").Value = '13.72'
$ws.Range("E46").Value = '  -3.00%  '
$ws.Range("E47").Value = '  -3.43%  '
$ws.Range("D48").Value = '0.5849'
$ws.Range("E48").Value = '  -4.60%  '
$ws.Range("D49").Value = '124.69'
$ws.Range("E49").Value = '  -4.35%  '
$ws.Range("D50").Value = '1.945'
$ws.Range("E50").Value = '  -4.38%  '
$ws.Range("D51").Value = '0.06850'
